# Revise analysis of chlorides: the cumulative-area (E) and cumulative-IC
# (F) running totals for row 4 (Upper Main Stem, S06B) previously only
# rolled forward row 3's cumulative total plus the current row's own
# Area_ac/IC_ac (C4/D4). The corrected analysis also needs to fold in
# row 2's (Blanchette Brook, S07) contribution, since S07 is a separate
# headwater subwatershed that drains into the same cumulative total at
# this point in the watershed, and was omitted from the running sum.
#
# Updating E4 and F4 automatically ripples through the dependent
# formulas in E5/F5 (next cumulative total) and H4/H5 (CumPctIC =
# CumIC_ac / CumArea_ac), all of which simply recalculate.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# CumArea_ac (E) and CumIC_ac (F) for row 4 now also add in row 2's
# Area_ac / IC_ac contribution.
$ws.Range("E4").Formula = "=E3+C4+C2"
$ws.Range("F4").Formula = "=F3+D4+D2"

# Leave the selection on H5, matching where the author's cursor ended up
# after reviewing the corrected CumPctIC figure.
$ws.Range("H5").Select()
